$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1895910780669145
$ws.Range("C2").Value = 0.5687732342007435
$ws.Range("J2").Value = 0.01858736059479554
$ws.Range("O2").Value = 0.003717472118959108
$ws.Range("P2").Value = 0.1263940520446097
$ws.Range("S2").Value = 0.09293680297397769

$ws.Range("B3").Value = 0.01265822784810127
$ws.Range("C3").Value = 0.02531645569620253
$ws.Range("J3").Value = 0.006329113924050633
$ws.Range("P3").Value = 0.7784810126582279
$ws.Range("S3").Value = 0.1772151898734177

$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.2708333333333333

$ws.Range("B6").Value = 0.04263565891472868
$ws.Range("D6").Value = 0.007751937984496124
$ws.Range("F6").Value = 0.05813953488372093
$ws.Range("J6").Value = 0.2790697674418605
$ws.Range("O6").Value = 0.003875968992248062
$ws.Range("Q6").Value = 0.1201550387596899
$ws.Range("R6").Value = 0.1124031007751938
$ws.Range("S6").Value = 0.375968992248062

$ws.Range("B7").Value = 0.09251101321585903
$ws.Range("D7").Value = 0.02643171806167401
$ws.Range("F7").Value = 0.04405286343612335
$ws.Range("J7").Value = 0.1013215859030837
$ws.Range("O7").Value = 0.02643171806167401
$ws.Range("Q7").Value = 0.1938325991189427
$ws.Range("R7").Value = 0.05286343612334802
$ws.Range("S7").Value = 0.4625550660792951

$ws.Range("B8").Value = 0.07792207792207792
$ws.Range("D8").Value = 0.00927643784786642
$ws.Range("F8").Value = 0.05565862708719851
$ws.Range("J8").Value = 0.1391465677179963
$ws.Range("O8").Value = 0.01298701298701299
$ws.Range("Q8").Value = 0.176252319109462
$ws.Range("R8").Value = 0.09276437847866419
$ws.Range("S8").Value = 0.4359925788497217

$ws.Range("B9").Value = 0.09714285714285714
$ws.Range("F9").Value = 0.1028571428571429
$ws.Range("J9").Value = 0.1371428571428571
$ws.Range("O9").Value = 0.01142857142857143
$ws.Range("Q9").Value = 0.1371428571428571
$ws.Range("R9").Value = 0.07428571428571429
$ws.Range("S9").Value = 0.44

$ws.Range("B10").Value = 0.08608321377331421
$ws.Range("D10").Value = 0.02582496413199426
$ws.Range("E10").Value = 0.001434720229555237
$ws.Range("F10").Value = 0.06671449067431851
$ws.Range("J10").Value = 0.1312769010043042
$ws.Range("O10").Value = 0.007890961262553802
$ws.Range("Q10").Value = 0.2087517934002869
$ws.Range("R10").Value = 0.08321377331420372
$ws.Range("S10").Value = 0.3888091822094691

$ws.Range("G11").Value = 0.1201117318435754
$ws.Range("J11").Value = 0.1005586592178771
$ws.Range("K11").Value = 0.1731843575418995
$ws.Range("L11").Value = 0.5782122905027933
$ws.Range("S11").Value = 0.02793296089385475

$ws.Range("G12").Value = 0.7570093457943925
$ws.Range("J12").Value = 0.1588785046728972
$ws.Range("K12").Value = 0.004672897196261682
$ws.Range("L12").Value = 0.01401869158878505
$ws.Range("S12").Value = 0.06542056074766354

$ws.Range("G13").Value = 0.7
$ws.Range("J13").Value = 0.28
$ws.Range("S13").Value = 0.02

$ws.Range("J14").Value = 1

$ws.Range("F15").Value = 0.02714932126696833
$ws.Range("H15").Value = 0.2126696832579185
$ws.Range("I15").Value = 0.06787330316742081
$ws.Range("J15").Value = 0.3529411764705883
$ws.Range("K15").Value = 0.06334841628959276
$ws.Range("M15").Value = 0.01809954751131222
$ws.Range("O15").Value = 0.06334841628959276
$ws.Range("S15").Value = 0.1945701357466063

$ws.Range("F16").Value = 0.02173913043478261
$ws.Range("H16").Value = 0.1684782608695652
$ws.Range("I16").Value = 0.06521739130434782
$ws.Range("J16").Value = 0.3858695652173913
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.03260869565217391
$ws.Range("N16").Value = 0.005434782608695652
$ws.Range("O16").Value = 0.07608695652173914
$ws.Range("S16").Value = 0.1576086956521739

$ws.Range("F17").Value = 0.02916666666666667
$ws.Range("H17").Value = 0.19375
$ws.Range("I17").Value = 0.07291666666666667
$ws.Range("J17").Value = 0.3979166666666666
$ws.Range("K17").Value = 0.1208333333333333
$ws.Range("M17").Value = 0.01875
$ws.Range("O17").Value = 0.04583333333333333
$ws.Range("S17").Value = 0.1208333333333333

$ws.Range("F18").Value = 0.02727272727272727
$ws.Range("H18").Value = 0.1954545454545455
$ws.Range("I18").Value = 0.08636363636363636
$ws.Range("J18").Value = 0.35
$ws.Range("K18").Value = 0.1363636363636364
$ws.Range("M18").Value = 0.01363636363636364
$ws.Range("O18").Value = 0.06363636363636363
$ws.Range("S18").Value = 0.1272727272727273

$ws.Range("F19").Value = 0.01875
$ws.Range("H19").Value = 0.2243055555555556
$ws.Range("I19").Value = 0.06736111111111111
$ws.Range("J19").Value = 0.3638888888888889
$ws.Range("K19").Value = 0.1215277777777778
$ws.Range("M19").Value = 0.01944444444444444
$ws.Range("N19").Value = 0.0006944444444444445
$ws.Range("O19").Value = 0.0763888888888889
$ws.Range("S19").Value = 0.1076388888888889
